# "cambios de mes de mayo" - roll the reporting period forward one quarter
# (Q4 2021 -> Q1 2022) on the single data row of the "Reporte de Formatos"
# sheet, and leave the view scrolled/selected where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ejercicio (year) 2021 -> 2022
$ws.Range("A8").Value = 2022

# Fecha de inicio del periodo que se informa: 2021-10-01 -> 2022-01-01
$ws.Range("B8").Value = 44562

# Fecha de término del periodo que se informa: 2021-12-31 -> 2022-03-31
$ws.Range("C8").Value = 44651

# Fecha de validación (F) and Fecha de actualización (G): 2022-01-10 -> 2022-04-08
$ws.Range("F8").Value = 44659
$ws.Range("G8").Value = 44659

# Move the view / selection to match where the author left the cursor
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("H8").Select()
